$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 303 (shifts existing rows 303:384 down to 304:385,
# Excel carries formatting down from the row above automatically).
$ws.Rows("303:303").Insert()

# Populate the newly inserted row 303 with the new record.
$ws.Range("A303").Value = 10
$ws.Range("B303").Value = "Vega Modelo de Temuco"
$ws.Range("C303").Value = "La Araucanía"
$ws.Range("D303").Value = 44754
$ws.Range("E303").Value = 9
$ws.Range("F303").Value = "Fruta"
$ws.Range("G303").Value = 100108
$ws.Range("H303").Value = "Tropicales y subtropicales"
$ws.Range("I303").Value = 100108002
$ws.Range("J303").Value = "Mango"
$ws.Range("K303").Value = "Sin especificar"
$ws.Range("L303").Value = "Primera"
$ws.Range("M303").Value = 155
$ws.Range("N303").Value = 8000
$ws.Range("O303").Value = 8000
$ws.Range("P303").Value = 8000
$ws.Range("Q303").Value = "$/bandeja 4 kilos"
$ws.Range("R303").Value = "Brasil"
$ws.Range("S303").Value = 2000
$ws.Range("T303").Value = 4
